$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 3773.2666
$ws.Range("J40").Value = 9499.75
$ws.Range("L40").Value = 9499.75
$ws.Range("N40").Value = -9849.75

# Row 41
$ws.Range("H41").Value = 893.3158
$ws.Range("I41").Value = 955.9231
$ws.Range("J41").Value = 757.6667
$ws.Range("K41").Value = 955.9231
$ws.Range("L41").Value = 757.6667
$ws.Range("M41").Value = -515.9231
$ws.Range("N41").Value = -1637.6667

# Row 62
$ws.Range("H62").Value = 1750
$ws.Range("I62").Value = 1500
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 1500
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -876
$ws.Range("N62").Value = -3248

# Row 65
$ws.Range("H65").Value = 1750
$ws.Range("I65").Value = 1500
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 7500
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -4380
$ws.Range("N65").Value = -16240

# Row 113
$ws.Range("H113").Value = 55758.5
$ws.Range("I113").Value = 75799.28999999999
$ws.Range("K113").Value = 75799.28999999999
$ws.Range("M113").Value = -72545.28999999999

# Row 116
$ws.Range("H116").Value = 7844.3706
$ws.Range("I116").Value = 7500
$ws.Range("K116").Value = 7500
$ws.Range("M116").Value = -4058

# Row 132
$ws.Range("H132").Value = 1518467.5
$ws.Range("I132").Value = 1518467.5
$ws.Range("K132").Value = 4555402.5
$ws.Range("M132").Value = -4552872.5

# Row 135
$ws.Range("H135").Value = 13158628
$ws.Range("I135").Value = 14706656
$ws.Range("K135").Value = 132359904
$ws.Range("M135").Value = -132357369

# Row 137
$ws.Range("H137").Value = 1317.1052
$ws.Range("I137").Value = 1298.4286
$ws.Range("J137").Value = 1369.4
$ws.Range("K137").Value = 3895.2858
$ws.Range("L137").Value = 4108.200000000001
$ws.Range("M137").Value = -1345.2858
$ws.Range("N137").Value = -9208.200000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1371.2
$ws.Range("I2").Value = 1344.1
$ws.Range("K2").Value = 1344.1
$ws.Range("M2").Value = -1231.1

# Row 45
$ws.Range("H45").Value = 2323.9443
$ws.Range("I45").Value = 1202.8462
$ws.Range("K45").Value = 1202.8462
$ws.Range("M45").Value = -825.8462

# Row 116
$ws.Range("H116").Value = 1371.2
$ws.Range("I116").Value = 1344.1
$ws.Range("K116").Value = 1344.1
$ws.Range("M116").Value = 949.9000000000001

# Row 132
$ws.Range("H132").Value = 3687.4688
$ws.Range("I132").Value = 3038.1924
$ws.Range("K132").Value = 9114.5772
$ws.Range("M132").Value = -6584.5772

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1371.2
$ws.Range("I3").Value = 1344.1
$ws.Range("K3").Value = 1344.1
$ws.Range("M3").Value = -1230.1

# Row 54
$ws.Range("H54").Value = 1145.25
$ws.Range("I54").Value = 1145.25
$ws.Range("K54").Value = 1145.25
$ws.Range("M54").Value = -661.25

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 279.94116
$ws.Range("I7").Value = 172.45454
$ws.Range("J7").Value = 477
$ws.Range("K7").Value = 172.45454
$ws.Range("L7").Value = 477
$ws.Range("M7").Value = -59.45454000000001
$ws.Range("N7").Value = -703

# Row 31
$ws.Range("H31").Value = 1590.76
$ws.Range("I31").Value = 1243.2
$ws.Range("J31").Value = 2112.1
$ws.Range("K31").Value = 1243.2
$ws.Range("L31").Value = 2112.1
$ws.Range("M31").Value = -948.2
$ws.Range("N31").Value = -2702.1

# Row 34
$ws.Range("H34").Value = 1590.76
$ws.Range("I34").Value = 1243.2
$ws.Range("J34").Value = 2112.1
$ws.Range("K34").Value = 1243.2
$ws.Range("L34").Value = 2112.1
$ws.Range("M34").Value = -1041.2
$ws.Range("N34").Value = -2516.1

# Row 58
$ws.Range("H58").Value = 43752810
$ws.Range("I58").Value = 16668868
$ws.Range("J58").Value = 125004620
$ws.Range("K58").Value = 16668868
$ws.Range("L58").Value = 125004620
$ws.Range("M58").Value = -16668665
$ws.Range("N58").Value = -125005026

# Row 99
$ws.Range("H99").Value = 4281.1665
$ws.Range("I99").Value = 4504.8
$ws.Range("J99").Value = 3163
$ws.Range("K99").Value = 4504.8
$ws.Range("L99").Value = 3163
$ws.Range("M99").Value = -3006.8
$ws.Range("N99").Value = -6159

# Row 107
$ws.Range("H107").Value = 1367.2122
$ws.Range("I107").Value = 1116.3334
$ws.Range("K107").Value = 1116.3334
$ws.Range("M107").Value = 803.6666

# Row 126
$ws.Range("H126").Value = 4281.1665
$ws.Range("I126").Value = 4504.8
$ws.Range("J126").Value = 3163
$ws.Range("K126").Value = 13514.4
$ws.Range("L126").Value = 9489
$ws.Range("M126").Value = -11044.4
$ws.Range("N126").Value = -14429

# Row 135
$ws.Range("H135").Value = 66274.586
$ws.Range("J135").Value = 66274.586
$ws.Range("L135").Value = 66274.586
$ws.Range("N135").Value = -76414.586

# Row 136
$ws.Range("H136").Value = 43752810
$ws.Range("I136").Value = 16668868
$ws.Range("J136").Value = 125004620
$ws.Range("K136").Value = 50006604
$ws.Range("L136").Value = 375013860
$ws.Range("M136").Value = -50004054
$ws.Range("N136").Value = -375018960

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 52.5
$ws.Range("J2").Value = 53
$ws.Range("L2").Value = 318
$ws.Range("N2").Value = -544

# Row 11
$ws.Range("H11").Value = 318974.1
$ws.Range("I11").Value = 846.5
$ws.Range("J11").Value = 3500250
$ws.Range("K11").Value = 2539.5
$ws.Range("L11").Value = 10500750
$ws.Range("M11").Value = -2399.5
$ws.Range("N11").Value = -10501030

# Row 34
$ws.Range("H34").Value = 3091
$ws.Range("I34").Value = 434
$ws.Range("J34").Value = 5748
$ws.Range("K34").Value = 1302
$ws.Range("L34").Value = 17244
$ws.Range("M34").Value = -1218
$ws.Range("N34").Value = -17412

# Row 39
$ws.Range("H39").Value = 6447
$ws.Range("J39").Value = 8167.8887
$ws.Range("L39").Value = 24503.6661
$ws.Range("N39").Value = -25091.6661

# Row 55
$ws.Range("H55").Value = 15252.223
$ws.Range("J55").Value = 19138.572
$ws.Range("L55").Value = 57415.716
$ws.Range("N55").Value = -57769.716

# Row 121
$ws.Range("H121").Value = 5558565
$ws.Range("J121").Value = 7693351.5
$ws.Range("L121").Value = 23080054.5
$ws.Range("N121").Value = -23082674.5

# Row 127
$ws.Range("H127").Value = 57176.375
$ws.Range("J127").Value = 57176.375
$ws.Range("L127").Value = 171529.125
$ws.Range("N127").Value = -181449.125

# Row 137
$ws.Range("H137").Value = 1962.7053
$ws.Range("I137").Value = 886
$ws.Range("K137").Value = 2658
$ws.Range("M137").Value = 2442

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2648.3
$ws.Range("I80").Value = 1989.5
$ws.Range("J80").Value = 3636.5
$ws.Range("K80").Value = 1989.5
$ws.Range("L80").Value = 3636.5
$ws.Range("M80").Value = -991.5
$ws.Range("N80").Value = -5632.5

# Row 83
$ws.Range("H83").Value = 2648.3
$ws.Range("I83").Value = 1989.5
$ws.Range("J83").Value = 3636.5
$ws.Range("K83").Value = 9947.5
$ws.Range("L83").Value = 18182.5
$ws.Range("M83").Value = -4955.5
$ws.Range("N83").Value = -28166.5

# Row 102
$ws.Range("H102").Value = 10316.5
$ws.Range("J102").Value = 11096
$ws.Range("L102").Value = 11096
$ws.Range("N102").Value = -14340

# Row 107
$ws.Range("H107").Value = 372.3846
$ws.Range("I107").Value = 372.3846
$ws.Range("K107").Value = 372.3846
$ws.Range("M107").Value = 1547.6154

# Row 122
$ws.Range("H122").Value = 88631.21000000001
$ws.Range("I122").Value = 115373.445
$ws.Range("J122").Value = 8404.5
$ws.Range("K122").Value = 346120.335
$ws.Range("L122").Value = 25213.5
$ws.Range("M122").Value = -343670.335
$ws.Range("N122").Value = -30113.5

# Row 132
$ws.Range("H132").Value = 1908825.8
$ws.Range("I132").Value = 2226130.8
$ws.Range("K132").Value = 6678392.399999999
$ws.Range("M132").Value = -6675862.399999999

$ws = $wb.Worksheets.Item("LTW")
# Row 13
$ws.Range("H13").Value = 7665.6665
$ws.Range("I13").Value = 496
$ws.Range("K13").Value = 496
$ws.Range("M13").Value = -356

# Row 22
$ws.Range("H22").Value = 4985.8
$ws.Range("J22").Value = 5701.143
$ws.Range("L22").Value = 5701.143
$ws.Range("N22").Value = -6291.143

# Row 27
$ws.Range("H27").Value = 4985.8
$ws.Range("J27").Value = 5701.143
$ws.Range("L27").Value = 5701.143
$ws.Range("N27").Value = -5915.143

# Row 40
$ws.Range("H40").Value = 5081.077
$ws.Range("I40").Value = 5717.5
$ws.Range("K40").Value = 5717.5
$ws.Range("M40").Value = -5581.5

# Row 55
$ws.Range("H55").Value = 2311.3
$ws.Range("I55").Value = 566.1111
$ws.Range("K55").Value = 566.1111
$ws.Range("M55").Value = -393.1111

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 4004
$ws.Range("I122").Value = 3589.6875
$ws.Range("K122").Value = 10769.0625
$ws.Range("M122").Value = -8319.0625

# Row 123
$ws.Range("H123").Value = 69238.22
$ws.Range("J123").Value = 69238.22
$ws.Range("L123").Value = 69238.22
$ws.Range("N123").Value = -79038.22

# Row 132
$ws.Range("H132").Value = 2944.1333
$ws.Range("I132").Value = 2909.7778
$ws.Range("K132").Value = 8729.3334
$ws.Range("M132").Value = -6199.3334

# Row 133
$ws.Range("H133").Value = 56396
$ws.Range("J133").Value = 56396
$ws.Range("L133").Value = 56396
$ws.Range("N133").Value = -66516
